$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold font, border, centered alignment) from G1 to the new H1 header cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New "Save" column: header label + data value
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
